$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")
$ws.Range("B2").Formula = "=""TRUE"""
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
